$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S24 Table")

# Update cell values per revised reconstruction-tools data
$ws.Range("B4").Value = [double]"59"
$ws.Range("E4").Value = [double]"518"
$ws.Range("F4").Value = [double]"722"
$ws.Range("G4").Value = [double]"762"
$ws.Range("E5").Value = [double]"0"
$ws.Range("F5").Value = [double]"10"
$ws.Range("G5").Value = [double]"15"
$ws.Range("E6").Value = [double]"2"
$ws.Range("F6").Value = [double]"3"
$ws.Range("G6").Value = [double]"3"
$ws.Range("B7").Value = [double]"4.704944178628389E-2"
$ws.Range("E7").Value = [double]"0.45438596491228073"
$ws.Range("F7").Value = [double]"0.48816768086544965"
$ws.Range("G7").Value = [double]"0.49544863459037713"
$ws.Range("B8").Value = [double]"0.35805422647527912"
$ws.Range("E8").Value = [double]"0"
$ws.Range("F8").Value = [double]"6.7613252197430695E-3"
$ws.Range("G8").Value = [double]"9.7529258777633299E-3"
$ws.Range("B9").Value = [double]"7.9744816586921851E-4"
$ws.Range("E9").Value = [double]"1.7543859649122807E-3"
$ws.Range("F9").Value = [double]"2.0283975659229209E-3"
$ws.Range("G9").Value = [double]"1.9505851755526658E-3"
$ws.Range("E12").Value = [double]"5"
$ws.Range("F12").Value = [double]"6"
$ws.Range("G12").Value = [double]"6"
$ws.Range("F13").Value = [double]"11"
$ws.Range("G13").Value = [double]"12"
$ws.Range("E14").Value = [double]"53.799999999997453"
$ws.Range("F14").Value = [double]"61.840000000000146"
$ws.Range("G14").Value = [double]"64"
$ws.Range("B19").Value = [double]"35.880000000000109"
$ws.Range("E19").Value = [double]"33.199999999999818"
$ws.Range("F19").Value = [double]"34.420000000000073"
$ws.Range("G19").Value = [double]"33.239999999999782"
$ws.Range("E22").Value = [double]"6"
$ws.Range("F22").Value = [double]"6"
$ws.Range("G22").Value = [double]"6"
$ws.Range("B23").Value = [double]"9"
$ws.Range("E23").Value = [double]"8"
$ws.Range("F23").Value = [double]"8"
$ws.Range("G23").Value = [double]"8"
$ws.Range("B24").Value = [double]"27.960000000000036"
$ws.Range("E24").Value = [double]"25.099999999999909"
$ws.Range("F24").Value = [double]"25.420000000000073"
$ws.Range("G24").Value = [double]"25.239999999999782"
$ws.Range("G28").Value = [double]"8"
$ws.Range("F29").Value = [double]"12"
$ws.Range("G29").Value = [double]"13"

# Restore the saved selection/active cell state
$ws.Range("E16").Select()
